$wb = $excel.ActiveWorkbook

# --- Sheet "simulations_scripts": add two new rows (14 and 15) ---
$wsScripts = $wb.Worksheets.Item("simulations_scripts")
$wsScripts.Range("A14").Value = "nrgacqweek"
$wsScripts.Range("A15").Value = "newnrgweek"
$wsScripts.Range("B14").Value = "energy acquired in that week (last day of week before - last day of the reference week)"
$wsScripts.Range("B15").Value = "energy renewed for that week (ccr * 10 of the first day of the reference week)"

# --- Sheet "Sheet1": remove column D content (D10 cell cleared) ---
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Range("D10").Clear()
$wsData.Range("C10").ClearFormats()

# --- Selections / active sheet state ---
[void]$wsData.Range("D13").Select()
[void]$wsScripts.Select()
[void]$wsScripts.Range("B16").Select()
